$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "download this file" note (F2) with the new URL/casing ---
$ws.Range("F2").Value = "You can download this file from https://njan-oru-malayali.com/Grocery-Expenses-in-Germany.xlsx?raw=true"

# --- Add "Price in Indian Rupee" (column E) conversion formulas for rows 2-27 ---
$ws.Range("E2").Formula = "=D2*76.88"
$ws.Range("E3:E27").Formula = "=D3*76.88"

# --- Row 27: Total row ---
$ws.Range("A27").Value = "Total"
$ws.Range("D27").Value = 72.94

# Keep E27's original (General) number format instead of the currency-style
# format that gets auto-applied when a formula referencing D27 is entered.
$ws.Range("A27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Column D width widened to fit the new "Total" values ---
$ws.Columns("D").ColumnWidth = 20.9296875

# --- Selection moved to F27 ---
$ws.Range("F27").Select() | Out-Null
